# Restore the "From" value of rule R30 (row 10) in the Rules sheet.
# C10 changes from 18 -> 1 (the cell's existing number format/style, s="20",
# is left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
